$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for rows 2-12 per repulled data
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = -5
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 4
$ws.Range("F11").Value = -4
$ws.Range("F12").Value = -4
